$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date/time-looking values to be stored as plain text, matching
# the inlineStr cells used throughout the rest of the sheet, instead of
# letting Excel auto-convert them into date/time serial numbers.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("B21").NumberFormat = "@"

$ws.Range("A21").Value = "2025-09-16"
$ws.Range("B21").Value = "15:22:21"
$ws.Range("C21").Value = "1.00 EUR = 1,722.0677"
